$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = -0.1338314788754218;  C = 2.04803740314168;   D = 20.30227337679758;  E = 4.505804409514197;  F = 4.605038435222934;  G = 23 }
    3  = @{ B = -0.2196048395615291;  C = 2.254442816470436;  D = 18.45943047502218;  E = 4.296443933652827;  F = 4.39180238609223;   G = 22 }
    4  = @{ B = -0.6722431955569657;  C = 1.734626294717863;  D = 10.27176916279258;  E = 3.204960087550636;  F = 3.211051410560245;  G = 21 }
    5  = @{ B = -0.2619961853625158;  C = 1.574840796137216;  D = 11.52157846430296;  E = 3.394345071483298;  F = 3.472135163382664;  G = 20 }
    6  = @{ B = -0.2304210235372018;  C = 1.666997890576107;  D = 10.06653556874298;  E = 3.172780416092955;  F = 3.251114276439574;  G = 19 }
    7  = @{ B = -0.2737229439002919;  C = 1.775059182383401;  D = 10.57130412137001;  E = 3.251354198079626;  F = 3.333738739228395;  G = 18 }
    8  = @{ B = -0.1641233777288165;  C = 1.781553791812374;  D = 11.51292075854603;  E = 3.393069518672736;  F = 3.493402108638973;  G = 17 }
    9  = @{ B = -0.1510345969195566;  C = 1.85455728114967;   D = 11.36885938571983;  E = 3.37177392268815;   F = 3.478857733797603;  G = 16 }
    10 = @{ B = -0.1166201009408896;  C = 1.99642669408382;   D = 12.99001923899349;  E = 3.604166927182131;  F = 3.728713727525796;  G = 15 }
    11 = @{ B = -0.07784734317373232; C = 1.990677678875613;  D = 13.16673245956295;  E = 3.628599242071649;  F = 3.764708707805368;  G = 14 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
